$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "twitter"
$ws.Range("M1").Value = "instagram"
$ws.Range("N1").Value = "twitch"
$ws.Range("O1").Value = "facebook"

$lmno = @(
    @(0,1,1,0),
    @(1,1,1,0),
    @(0,1,0,1),
    @(1,1,0,1),
    @(1,1,1,0),
    @(1,1,1,0),
    @(1,1,0,1),
    @(1,0,1,0),
    @(0,1,0,1),
    @(0,1,0,0),
    @(0,0,0,1),
    @(0,1,0,0),
    @(1,0,0,1),
    @(1,0,0,1),
    @(0,0,0,0),
    @(1,1,0,0),
    @(1,1,0,1),
    @(1,1,0,1),
    @(0,1,0,1),
    @(1,0,0,1),
    @(1,1,0,0),
    @(1,0,0,0),
    @(0,0,0,0),
    @(1,0,0,1),
    @(0,0,0,0),
    @(1,1,0,1),
    @(0,0,0,1),
    @(0,1,0,0),
    @(0,1,0,0),
    @(0,1,0,1),
    @(1,1,1,1),
    @(1,1,0,1),
    @(1,1,0,0),
    @(0,0,1,1),
    @(0,0,0,0),
    @(1,1,0,0),
    @(1,0,0,0),
    @(1,1,0,1),
    @(1,1,0,1),
    @(1,0,0,1),
    @(1,1,0,1),
    @(0,0,0,0),
    @(1,1,0,1),
    @(1,1,0,1),
    @(1,1,0,1),
    @(1,1,0,1),
    @(0,0,0,0),
    @(1,1,0,1),
    @(1,1,0,1),
    @(1,1,0,0)
)

$startRow = 2
for ($i = 0; $i -lt $lmno.Count; $i++) {
    $r = $startRow + $i
    $vals = $lmno[$i]
    $ws.Cells.Item($r, 12).Value = $vals[0]
    $ws.Cells.Item($r, 13).Value = $vals[1]
    $ws.Cells.Item($r, 14).Value = $vals[2]
    $ws.Cells.Item($r, 15).Value = $vals[3]
}

$null = $ws.Range("F24").Select()
